$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0.001
$ws.Range("K14").Value = 465
$ws.Range("L14").Value = 0.0009300000000000001
